$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the publication-name text in column B (dropped " - Ovid search - ") ---
$ws.Range("B2").Value = "ICER RRMM 2022 report - ICER - 4/11/2022"
$ws.Range("B3").Value = "IC AML - Pfizer - 5/10/2021"
$ws.Range("B4").Value = "ICER RRMM 2022 report - ICER - 4/11/2022"
$ws.Range("B5").Value = "ICER RRMM 2022 report - ICER - 4/11/2022"
$ws.Range("B6").Value = "ICER RRMM 2022 report - ICER - 4/11/2022"
$ws.Range("B7").Value = "ICER RRMM 2022 report - ICER - 4/11/2022"

# --- Remember the header row's existing format so it can be restored below ---
$headerStyle = $ws.Range("A1").Style

# --- Strip the (unused/no-op) explicit style from columns A:D, which also clears
#     the per-cell style that used to sit on B2:D7 ---
$ws.Columns("A:D").ClearFormats()

# --- Restore the header row (row 1) formatting that ClearFormats above removed ---
$ws.Range("A1:D1").Style = $headerStyle

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()
